$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data range (A1:E11) first so stale cells (E column, rows 3-11) are removed
$ws.Range("A1:E11").Clear()

# Header row
$ws.Range("A1").Value = "Qtd_Nós"
$ws.Range("B1").Value = "Ativos"
$ws.Range("C1").Value = "Distancia"
$ws.Range("D1").Value = "Tempo"

# Data row
$ws.Range("A2").Value = 81
$ws.Range("B2").Value = 35
$ws.Range("C2").Value = 11538
$ws.Range("D2").Value = 0.2142186164855957
